$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental = "true" (as literal text, matching existing "Case Sensitive" true value)
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Date updated
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"

# Compositional = "false" (as literal text)
$ws.Range("B18").Formula = "=""false"""
$ws.Range("B18").Copy()
$ws.Range("B18").PasteSpecial(-4163)

$excel.CutCopyMode = $false
